# Generate Report for Handback
# Adds a new handback record (a2d5025a-aab5-4faf-b07f-34fc40b4d579) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$fileGuid = "a2d5025a-aab5-4faf-b07f-34fc40b4d579"
$xlfHash  = "bbfcbdecefe237f66a80da139c581d554f5ab7f9"

$mdName       = "$fileGuid.md"
$zhcnXlfName  = "$fileGuid.$xlfHash.zh-cn.xlf"
$dedeXlfName  = "$fileGuid.$xlfHash.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$includeText  = "Include"

$zhcnHandoffDt   = "2016-03-03 10:12:41"
$zhcnHandbackDt  = "2016-03-03 10:13:33"
$dedeHandoffDt   = "2016-03-03 10:12:53"
$dedeHandbackDt  = "2016-03-03 10:13:56"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/843e6c2432bb1e1e900f8fb29f45b4231d538805/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/843e6c2432bb1e1e900f8fb29f45b4231d538805/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZhCn.Range("B4").Value = $statusInSync

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b8a873a9c1c48e21ea9d223a9e03165926094e9c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhcnXlfName",
    "",
    "",
    $zhcnXlfName
) | Out-Null

$wsZhCn.Range("D4").Value = $zhcnHandoffDt

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0d64317afb368cbf170c8945a6e2e176ed83eeb7/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc37e51566e7e2300bfdd16d9679b1e4c0b0db27/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhcnXlfName",
    "",
    "",
    $zhcnXlfName
) | Out-Null

$wsZhCn.Range("G4").Value = $zhcnHandbackDt
$wsZhCn.Range("H4").Value = $includeText

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/843e6c2432bb1e1e900f8fb29f45b4231d538805/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDeDe.Range("B4").Value = $statusInSync

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5aa6176684ea45b34a0430fbb6cd1170556a3be7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$dedeXlfName",
    "",
    "",
    $dedeXlfName
) | Out-Null

$wsDeDe.Range("D4").Value = $dedeHandoffDt

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fb5a3842744c2d1bf6dee4cc3dd716739477defc/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f7f8ad395035860db5d1970aae539909d7c67d3b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$dedeXlfName",
    "",
    "",
    $dedeXlfName
) | Out-Null

$wsDeDe.Range("G4").Value = $dedeHandbackDt
$wsDeDe.Range("H4").Value = $includeText

Write-Output "Handback row added for $fileGuid"
